$p = $ppt.ActivePresentation

# --- Slide 7: "Learn Status Continued" table (by department) ---
$s7 = $p.Slides.Item(7)
$tbl7 = $s7.Shapes.Item(3).Table

$tbl7.Cell(2, 3).Shape.TextFrame.TextRange.Text = "172"    # Arcolab: Not on Learn 171 -> 172
$tbl7.Cell(4, 3).Shape.TextFrame.TextRange.Text = "338"    # Alathur: Not on Learn 214 -> 338
$tbl7.Cell(5, 3).Shape.TextFrame.TextRange.Text = "236"    # Corporate: Not on Learn 160 -> 236
$tbl7.Cell(6, 3).Shape.TextFrame.TextRange.Text = "67"     # Other: Not on Learn 64 -> 67
$tbl7.Cell(8, 3).Shape.TextFrame.TextRange.Text = "211"    # R&d: Not on Learn 209 -> 211
$tbl7.Cell(9, 3).Shape.TextFrame.TextRange.Text = "139"    # Seml: Not on Learn 137 -> 139
$tbl7.Cell(10, 3).Shape.TextFrame.TextRange.Text = "3190"  # Total: Not on Learn 2680 -> 3190
$tbl7.Cell(10, 4).Shape.TextFrame.TextRange.Text = "3192"  # Total: Total 2890 -> 3192

# --- Slide 11: "Learn Status Continued" table (subset) ---
$s11 = $p.Slides.Item(11)
$tbl11 = $s11.Shapes.Item(3).Table

$tbl11.Cell(2, 3).Shape.TextFrame.TextRange.Text = "286"   # Ucl: Not on Learn 282 -> 286
$tbl11.Cell(4, 3).Shape.TextFrame.TextRange.Text = "3192"  # Total: Not on Learn 298 -> 3192
$tbl11.Cell(4, 4).Shape.TextFrame.TextRange.Text = "3192"  # Total: Total 302 -> 3192
